$d = $word.ActiveDocument

$replacements = @(
    @{old = "680÷9=75, 5"; new = "999÷7=142, 5"},
    @{old = "170÷3=56, 2"; new = "508÷6=84, 4"},
    @{old = "923÷7=131, 6"; new = "186÷3=62, 0"},
    @{old = "980÷8=122, 4"; new = "428÷5=85, 3"},
    @{old = "114÷6=19, 0"; new = "323÷7=46, 1"},
    @{old = "300÷9=33, 3"; new = "473÷8=59, 1"},
    @{old = "517÷7=73, 6"; new = "670÷6=111, 4"},
    @{old = "332÷5=66, 2"; new = "261÷8=32, 5"},
    @{old = "244÷8=30, 4"; new = "323÷8=40, 3"},
    @{old = "463÷4=115, 3"; new = "116÷5=23, 1"},
    @{old = "587÷7=83, 6"; new = "636÷2=318, 0"},
    @{old = "370÷9=41, 1"; new = "566÷4=141, 2"},
    @{old = "444÷6=74, 0"; new = "837÷2=418, 1"},
    @{old = "273÷2=136, 1"; new = "415÷8=51, 7"},
    @{old = "372÷9=41, 3"; new = "294÷9=32, 6"},
    @{old = "593÷5=118, 3"; new = "132÷2=66, 0"},
    @{old = "592÷2=296, 0"; new = "307÷4=76, 3"},
    @{old = "799÷9=88, 7"; new = "728÷5=145, 3"},
    @{old = "935÷5=187, 0"; new = "306÷3=102, 0"},
    @{old = "742÷7=106, 0"; new = "176÷2=88, 0"},
    @{old = "510÷4=127, 2"; new = "625÷2=312, 1"},
    @{old = "724÷5=144, 4"; new = "644÷7=92, 0"},
    @{old = "361÷2=180, 1"; new = "805÷8=100, 5"},
    @{old = "906÷9=100, 6"; new = "558÷8=69, 6"},
    @{old = "226÷9=25, 1"; new = "902÷2=451, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
